$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("F1").Style = $ws.Range("E1").Style

$ws.Range("F2").Value = "2021-10-05 13:38:30.922050"
$ws.Range("F3").Value = "2021-10-05 13:38:30.922060"
$ws.Range("F4").Value = "2021-10-05 13:38:30.922063"
$ws.Range("F5").Value = "2021-10-05 13:38:30.922066"
$ws.Range("F6").Value = "2021-10-05 13:38:30.922068"
$ws.Range("F7").Value = "2021-10-05 13:38:30.922071"
$ws.Range("F8").Value = "2021-10-05 13:38:30.922073"
